# Update the three most-recent weekly rows (395-397) with the new week's
# data, then insert three rows to make room for what used to be there,
# re-entering the values that are being "pushed down" one week so the
# historical series keeps growing at the bottom (rows 428-430 are new).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Overwrite rows 395-397 with the new week's figures -----------------
# Row 395: Especial
$ws.Cells.Item(395, 4).Value = 44776
$ws.Cells.Item(395, 13).Value = 25
$ws.Cells.Item(395, 14).Value = 30000
$ws.Cells.Item(395, 15).Value = 30000
$ws.Cells.Item(395, 16).Value = 30000
$ws.Cells.Item(395, 19).Value = 1667

# Row 396: Primera
$ws.Cells.Item(396, 4).Value = 44776
$ws.Cells.Item(396, 13).Value = 50
$ws.Cells.Item(396, 14).Value = 25000
$ws.Cells.Item(396, 15).Value = 25000
$ws.Cells.Item(396, 16).Value = 25000
$ws.Cells.Item(396, 19).Value = 1389

# Row 397: Segunda
$ws.Cells.Item(397, 4).Value = 44776
$ws.Cells.Item(397, 13).Value = 40
$ws.Cells.Item(397, 14).Value = 18000
$ws.Cells.Item(397, 15).Value = 18000
$ws.Cells.Item(397, 16).Value = 18000
$ws.Cells.Item(397, 19).Value = 1000

# --- 2. Make room: insert 3 blank rows at row 398 ---------------------------
# This pushes the former rows 398-427 down to 401-430 intact.
$ws.Rows.Item(398).Insert()
$ws.Rows.Item(398).Insert()
$ws.Rows.Item(398).Insert()

# --- 3. Re-populate the newly inserted rows 398-400 -------------------------
# (same values the old rows 398-400 had, which now also live at 401-403)

# Row 398: Especial / Provincia de Melipilla
$ws.Cells.Item(398, 1).Value = 6
$ws.Cells.Item(398, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(398, 3).Value = "Metropolitana"
$ws.Cells.Item(398, 4).Value = 44494
$ws.Cells.Item(398, 5).Value = 13
$ws.Cells.Item(398, 6).Value = "Fruta"
$ws.Cells.Item(398, 7).Value = 100107
$ws.Cells.Item(398, 8).Value = "Otros"
$ws.Cells.Item(398, 9).Value = 100107011
$ws.Cells.Item(398, 10).Value = "Tuna"
$ws.Cells.Item(398, 11).Value = "Sin especificar"
$ws.Cells.Item(398, 12).Value = "Especial"
$ws.Cells.Item(398, 13).Value = 125
$ws.Cells.Item(398, 14).Value = 32000
$ws.Cells.Item(398, 15).Value = 32000
$ws.Cells.Item(398, 16).Value = 32000
$ws.Cells.Item(398, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(398, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(398, 19).Value = 1778
$ws.Cells.Item(398, 20).Value = 18

# Row 399: Primera / Provincia de Melipilla
$ws.Cells.Item(399, 1).Value = 6
$ws.Cells.Item(399, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(399, 3).Value = "Metropolitana"
$ws.Cells.Item(399, 4).Value = 44494
$ws.Cells.Item(399, 5).Value = 13
$ws.Cells.Item(399, 6).Value = "Fruta"
$ws.Cells.Item(399, 7).Value = 100107
$ws.Cells.Item(399, 8).Value = "Otros"
$ws.Cells.Item(399, 9).Value = 100107011
$ws.Cells.Item(399, 10).Value = "Tuna"
$ws.Cells.Item(399, 11).Value = "Sin especificar"
$ws.Cells.Item(399, 12).Value = "Primera"
$ws.Cells.Item(399, 13).Value = 20
$ws.Cells.Item(399, 14).Value = 20000
$ws.Cells.Item(399, 15).Value = 20000
$ws.Cells.Item(399, 16).Value = 20000
$ws.Cells.Item(399, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(399, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(399, 19).Value = 1111
$ws.Cells.Item(399, 20).Value = 18

# Row 400: Segunda / Provincia de Melipilla
$ws.Cells.Item(400, 1).Value = 6
$ws.Cells.Item(400, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(400, 3).Value = "Metropolitana"
$ws.Cells.Item(400, 4).Value = 44494
$ws.Cells.Item(400, 5).Value = 13
$ws.Cells.Item(400, 6).Value = "Fruta"
$ws.Cells.Item(400, 7).Value = 100107
$ws.Cells.Item(400, 8).Value = "Otros"
$ws.Cells.Item(400, 9).Value = 100107011
$ws.Cells.Item(400, 10).Value = "Tuna"
$ws.Cells.Item(400, 11).Value = "Sin especificar"
$ws.Cells.Item(400, 12).Value = "Segunda"
$ws.Cells.Item(400, 13).Value = 5
$ws.Cells.Item(400, 14).Value = 15000
$ws.Cells.Item(400, 15).Value = 15000
$ws.Cells.Item(400, 16).Value = 15000
$ws.Cells.Item(400, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(400, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(400, 19).Value = 833
$ws.Cells.Item(400, 20).Value = 18
